$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 185-186; this shifts the existing rows 185-253
# down to 187-255, preserving their values/styles.
$ws.Rows("185:186").Insert()

# New row 185: Camote, 1a (cosecha), new sampling date.
$ws.Range("A185").Value = 7
$ws.Range("B185").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C185").Value = 'Ñuble'
$ws.Range("D185").Value = 45006
$ws.Range("E185").Value = 16
$ws.Range("F185").Value = 100112045
$ws.Range("G185").Value = 'Zapallo'
$ws.Range("H185").Value = 'Camote'
$ws.Range("I185").Value = '1a (cosecha)'
$ws.Range("J185").Value = 250
$ws.Range("K185").Value = 500
$ws.Range("L185").Value = 500
$ws.Range("M185").Value = 500
$ws.Range("N185").Value = '$/kilo (volumen en unidades)'
$ws.Range("O185").Value = 'Región del Maule'
$ws.Range("P185").Value = 500
$ws.Range("Q185").Value = 1
$ws.Range("R185").Value = 'Hortaliza'

# New row 186: Paine, 1a (cosecha), same new sampling date.
$ws.Range("A186").Value = 7
$ws.Range("B186").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C186").Value = 'Ñuble'
$ws.Range("D186").Value = 45006
$ws.Range("E186").Value = 16
$ws.Range("F186").Value = 100112045
$ws.Range("G186").Value = 'Zapallo'
$ws.Range("H186").Value = 'Paine'
$ws.Range("I186").Value = '1a (cosecha)'
$ws.Range("J186").Value = 400
$ws.Range("K186").Value = 350
$ws.Range("L186").Value = 350
$ws.Range("M186").Value = 350
$ws.Range("N186").Value = '$/kilo (volumen en unidades)'
$ws.Range("O186").Value = 'Región del Maule'
$ws.Range("P186").Value = 350
$ws.Range("Q186").Value = 1
$ws.Range("R186").Value = 'Hortaliza'
